# Add "2022-Q4" quarter data:
#  1. Insert a new summary row at the top of the "总计" (total) sheet's data
#     for "2022-Q4" (26 funds, 1.72 billion yuan held), shifting existing
#     rows down by one.
#  2. Insert a brand-new worksheet named "2022-Q4" right after "总计"
#     holding the per-fund holdings table for that quarter.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1. "总计" sheet: insert new row 2 with the 2022-Q4 summary figures.
# ---------------------------------------------------------------------
$total = $wb.Worksheets.Item("总计")
$total.Rows.Item(2).Insert()

# The inserted row copies formatting from row 1 (header) by default.
# Restore the look used by the rest of the table: B2:D2 should be the
# plain (unstyled) data cells used elsewhere in those columns, and A2
# should match the bold/centered/bordered look already used by A3:A9
# (copied straight from a sibling cell so the same style id is reused).
$total.Range("B2:D2").Style = "Normal"
$total.Range("A3").Copy()
$total.Range("A2").PasteSpecial(-4122)

$total.Range("B2").Value = "2022-Q4"
$total.Range("C2").Value = 26
$total.Range("D2").Value = 1.72

# Renumber the running index in column A (0,1,2,...) for every data row
# now that an extra row sits at the top.
for ($r = 2; $r -le 9; $r++) {
    $total.Range("A" + $r).Value = $r - 2
}

# ---------------------------------------------------------------------
# 2. New "2022-Q4" sheet, inserted right after "总计".
# ---------------------------------------------------------------------
$q4 = $wb.Worksheets.Add($null, $total)
$q4.Name = "2022-Q4"

$headers = @("基金代码", "基金名称", "基金规模", "股票总仓位", "仓位占比", "持有市值(亿元)", "仓位排名")
for ($i = 0; $i -lt $headers.Length; $i++) {
    $cell = $q4.Cells.Item(1, $i + 2)
    $cell.Value = $headers[$i]
}
$headerRange = $q4.Range("B1:H1")
$headerRange.Font.Bold = $true
$headerRange.HorizontalAlignment = -4108
$headerRange.VerticalAlignment = -4160
$headerRange.Borders.LineStyle = 1

# Ensure the fund-code / text-looking numeric columns stay text (so
# leading zeros in fund codes like "004374" survive) before the values
# are written.
$q4.Range("B2:G27").NumberFormat = "@"

$rows = @(
    @("004374", "华泰保兴吉年丰混合A",               "5.90",  "94.77", "5.61", "0.3310", 2),
    @("006642", "华泰保兴吉年利定期开放混合",         "6.73",  "93.27", "4.71", "0.3170", 2),
    @("005313", "万家中证1000指数增强A",              "22.07", "94.13", "0.99", "0.2185", 7),
    @("005314", "万家中证1000指数增强C",              "19.61", "94.13", "0.99", "0.1941", 7),
    @("005904", "华泰保兴成长优选混合A",              "3.52",  "70.92", "3.64", "0.1281", 5),
    @("004375", "华泰保兴吉年丰混合C",                "1.34",  "94.77", "5.61", "0.0752", 2),
    @("015963", "汇安品质优选混合A",                  "2.48",  "86.98", "2.81", "0.0697", 10),
    @("001320", "工银丰盈回报灵活配置混合A",          "1.31",  "89.93", "4.91", "0.0643", 10),
    @("002212", "嘉实新起航灵活配置混合A",            "1.05",  "79.36", "5.86", "0.0615", 3),
    @("004050", "华夏新锦升灵活配置混合A",            "1.35",  "65.01", "4.37", "0.0590", 5),
    @("015635", "汇安价值先锋混合A",                  "1.41",  "84.53", "2.99", "0.0422", 10),
    @("014999", "华泰保兴吉年盈混合A",                "0.84",  "84.95", "4.19", "0.0352", 4),
    @("180028", "银华永祥灵活配置混合",                "0.70",  "77.51", "5.02", "0.0351", 5),
    @("013347", "工银丰盈回报灵活配置混合C",          "0.67",  "89.93", "4.91", "0.0329", 10),
    @("015964", "汇安品质优选混合C",                  "0.70",  "86.98", "2.81", "0.0197", 10),
    @("002159", "东吴国企改革主题灵活配置混合A",      "0.19",  "91.72", "5.95", "0.0113", 9),
    @("012615", "东吴国企改革主题灵活配置混合C",      "0.13",  "91.72", "5.95", "0.0077", 9),
    @("005000", "泰康泉林量化价值精选混合A",          "0.31",  "89.21", "1.61", "0.0050", 8),
    @("005905", "华泰保兴成长优选混合C",              "0.11",  "70.92", "3.64", "0.0040", 5),
    @("015636", "汇安价值先锋混合C",                  "0.09",  "84.53", "2.99", "0.0027", 10),
    @("005111", "泰康泉林量化价值精选混合C",          "0.14",  "89.21", "1.61", "0.0023", 8),
    @("004051", "华夏新锦升灵活配置混合C",            "0.05",  "65.01", "4.37", "0.0022", 5),
    @("007315", "汇安嘉盈一年持有期债券A",            "0.16",  "24.33", "1.15", "0.0018", 8),
    @("010270", "汇安嘉盈一年持有期债券C",            "0.15",  "24.33", "1.15", "0.0017", 8),
    @("015000", "华泰保兴吉年盈混合C",                "0.02",  "84.95", "4.19", "0.0008", 4),
    @("016264", "嘉实新起航灵活配置混合C",            "0.01",  "79.36", "5.86", "0.0006", 3)
)

for ($i = 0; $i -lt $rows.Length; $i++) {
    $r = $i + 2
    $row = $rows[$i]

    $q4.Cells.Item($r, 1).Value = $i
    $q4.Cells.Item($r, 2).Value = $row[0]
    $q4.Cells.Item($r, 3).Value = $row[1]
    $q4.Cells.Item($r, 4).Value = $row[2]
    $q4.Cells.Item($r, 5).Value = $row[3]
    $q4.Cells.Item($r, 6).Value = $row[4]
    $q4.Cells.Item($r, 7).Value = $row[5]
    $q4.Cells.Item($r, 8).Value = $row[6]
}

$idxRange = $q4.Range("A2:A27")
$idxRange.Font.Bold = $true
$idxRange.HorizontalAlignment = -4108
$idxRange.VerticalAlignment = -4160
$idxRange.Borders.LineStyle = 1

# Keep the originally-active "总计" tab selected (adding the sheet would
# otherwise leave the brand-new "2022-Q4" tab focused).
$total.Activate()
[void]$total.Range("A1").Select()

